# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.404.72'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.227.64'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.85'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.46'
$ws.Range('E6').Value = '  -4.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.563'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.494'
$ws.Range('E9').Value = '  -6.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.15'
$ws.Range('E10').Value = '  -6.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0785'
$ws.Range('E11').Value = '  -3.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.97'
$ws.Range('E12').Value = '  -4.31%  '
$ws.Range('D14').Value = '2.568.89'
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').Value = '2.226.36'
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.46'
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.781'
$ws.Range('E17').Value = '  -7.01%  '
$ws.Range('D18').Value = '44.400.39'
$ws.Range('E18').Value = '  +0.41%  '
$ws.Range('D19').Value = '0.0₃0914'
$ws.Range('E19').Value = '  -5.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.94'
$ws.Range('E20').Value = '  -7.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  -9.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.74'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.74'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.81'
$ws.Range('E24').Value = '  -6.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  -7.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.24'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.42'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.38'
$ws.Range('E29').Value = '  -5.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.53'
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '148.56'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.41'
$ws.Range('E32').Value = '  -9.63%  '
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0752'
$ws.Range('E34').Value = '  -6.05%  '
$ws.Range('E35').Value = '  -3.71%  '
$ws.Range('E36').Value = '  -10.69%  '
$ws.Range('E37').Value = '  -5.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.69'
$ws.Range('E38').Value = '  -6.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0304'
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.21'
$ws.Range('E40').Value = '  -6.94%  '
$ws.Range('E41').Value = '  -7.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.27'
$ws.Range('E42').Value = '  -9.12%  '
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = '1.815.16'
$ws.Range('E44').Value = '  +3.69%  '
$ws.Range('E45').Value = '  +11.87%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '14.20'
$ws.Range('E46').Value = '  +9.39%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.180'
$ws.Range('E47').Value = '  -7.45%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '74.69'
$ws.Range('E48').Value = '  -7.70%  '
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.84'
$ws.Range('E49').Value = '  +11.46%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.37'
$ws.Range('E50').Value = '  -5.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.00'
$ws.Range('E51').Value = '  -5.86%  '
